$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = 1.39214
$ws.Range("E4").Value = 0.04983
$ws.Range("D5").Value = 1.39214
$ws.Range("E5").Value = 0.04191
$ws.Range("D6").Value = 1.39214
$ws.Range("E6").Value = 0.04041
$ws.Range("D7").Value = 1.39214
$ws.Range("E7").Value = 0.04139
$ws.Range("D8").Value = 1.39214
$ws.Range("E8").Value = 0.04125
$ws.Range("D9").Value = 1.39214
$ws.Range("E9").Value = 0.0438
$ws.Range("D11").Value = 3.74431
$ws.Range("E11").Value = 0.08907
$ws.Range("D12").Value = 3.74431
$ws.Range("E12").Value = 0.07274
$ws.Range("D13").Value = 3.74431
$ws.Range("E13").Value = 0.07394
$ws.Range("D14").Value = 3.74431
$ws.Range("E14").Value = 0.07341
$ws.Range("D15").Value = 3.74431
$ws.Range("E15").Value = 0.07361
$ws.Range("D16").Value = 3.74431
$ws.Range("E16").Value = 0.07516
$ws.Range("D18").Value = 6.04087
$ws.Range("E18").Value = 0.13596
$ws.Range("D19").Value = 6.04087
$ws.Range("E19").Value = 0.10761
$ws.Range("D20").Value = 6.04087
$ws.Range("E20").Value = 0.10795
$ws.Range("D21").Value = 6.04087
$ws.Range("E21").Value = 0.10791
$ws.Range("D22").Value = 6.04087
$ws.Range("E22").Value = 0.11034
$ws.Range("D23").Value = 6.04087
$ws.Range("E23").Value = 0.11333
$ws.Range("D25").Value = 7.54378
$ws.Range("E25").Value = 0.18132
$ws.Range("D26").Value = 7.54378
$ws.Range("E26").Value = 0.14308
$ws.Range("D27").Value = 7.54378
$ws.Range("E27").Value = 0.14356
$ws.Range("D28").Value = 7.54378
$ws.Range("E28").Value = 0.14406
$ws.Range("D29").Value = 7.54378
$ws.Range("E29").Value = 0.14661
$ws.Range("D30").Value = 7.54378
$ws.Range("E30").Value = 0.14924
$ws.Range("D35").Value = 0.70481
$ws.Range("E35").Value = 0.05393
$ws.Range("D36").Value = 0.70481
$ws.Range("E36").Value = 0.04867
$ws.Range("D37").Value = 0.70481
$ws.Range("E37").Value = 0.04292
$ws.Range("D38").Value = 0.70481
$ws.Range("E38").Value = 0.04187
$ws.Range("D39").Value = 0.70481
$ws.Range("E39").Value = 0.04493
$ws.Range("D40").Value = 0.70481
$ws.Range("E40").Value = 0.05469
$ws.Range("D42").Value = 1.56794
$ws.Range("E42").Value = 0.08151
$ws.Range("D43").Value = 1.56794
$ws.Range("E43").Value = 0.06928
$ws.Range("D44").Value = 1.56794
$ws.Range("E44").Value = 0.06246
$ws.Range("D45").Value = 1.56794
$ws.Range("E45").Value = 0.06106
$ws.Range("D46").Value = 1.56794
$ws.Range("E46").Value = 0.06298
$ws.Range("D47").Value = 1.56794
$ws.Range("E47").Value = 0.09201
$ws.Range("D49").Value = 2.35501
$ws.Range("E49").Value = 0.10484
$ws.Range("D50").Value = 2.35501
$ws.Range("E50").Value = 0.08806
$ws.Range("D51").Value = 2.35501
$ws.Range("E51").Value = 0.07792
$ws.Range("D52").Value = 2.35501
$ws.Range("E52").Value = 0.07653
$ws.Range("D53").Value = 2.35501
$ws.Range("E53").Value = 0.08272
$ws.Range("D54").Value = 2.35501
$ws.Range("E54").Value = 0.12428
$ws.Range("D56").Value = 3.1847
$ws.Range("E56").Value = 0.12624
$ws.Range("D57").Value = 3.1847
$ws.Range("E57").Value = 0.10406
$ws.Range("D58").Value = 3.1847
$ws.Range("E58").Value = 0.09318
$ws.Range("D59").Value = 3.1847
$ws.Range("E59").Value = 0.09128
$ws.Range("D60").Value = 3.1847
$ws.Range("E60").Value = 0.09376
$ws.Range("D61").Value = 3.1847
$ws.Range("E61").Value = 0.15326

# Reset the scroll position (sheet view no longer pinned to A37)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
